# Generate Report for Handoff
#
# A fresh localization-status report was generated: the files that were
# handed off together (b8e4201d, 29d10456, 4eb5ff2b, 9b433dc5, a6aeea7e,
# dd20e269, f583b714, fa066fdc - i.e. table rows 7, 10, 11, 12, 13, 14, 15,
# 16) now all carry a single, newer "Latest Handoff Datetime" (column D)
# for each localized-language sheet.

$wb = $excel.ActiveWorkbook

$rows = @(7, 10, 11, 12, 13, 14, 15, 16)

$ws_zhcn = $wb.Worksheets.Item("zh-cn")
foreach ($r in $rows) {
    $ws_zhcn.Range("D$r").Value = "2016-03-10 22:29:43"
}

$ws_dede = $wb.Worksheets.Item("de-de")
foreach ($r in $rows) {
    $ws_dede.Range("D$r").Value = "2016-03-10 22:29:50"
}
